# Update the dSF (column F) values for several rows to reflect the
# repulled data / recalculated mean, per the commit message
# "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    7  = 0
    8  = -6
    10 = -3
    12 = 1
    15 = -7
    16 = 4
    21 = -2
    23 = -1
    27 = -11
    28 = -1
    31 = 0
    38 = 0
    43 = -1
    46 = 0
    48 = 2
    51 = 4
    55 = -1
    63 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
